# Actualización automática 2025-11-10 16:30:09
#
# Updates sales figures for GUERRERO FAREZ FABIAN MAURICIO: a new
# "noviembre" (PORCELANATO) sale of 149.69 for client "AGUIMPORT-AGUILAR
# IMPORTACIONES S.A.S." and an increase from 33.7 to 56.86 for client
# "ORTEGA ROMAN KLEBER ERWIN", plus the downstream totals / compliance
# figures that these roll up into.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (sales by product group, per client)
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# PORCELANATO column (M) for the two affected clients
$wsGrupo.Range("M5").Value  = 149.69   # AGUIMPORT-AGUILAR IMPORTACIONES S.A.S.
$wsGrupo.Range("M36").Value = 56.86    # ORTEGA ROMAN KLEBER ERWIN

# Count of clients with PORCELANATO sales (out of 54) increases by one
$wsGrupo.Range("M56").Value = "7 de 54"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (sales by month, per client)
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# "noviembre" column (F) for the same two clients
$wsMensual.Range("F5").Value  = 149.69
$wsMensual.Range("F36").Value = 56.86

# Monthly total row
$wsMensual.Range("F60").Value = 10150.5

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" (budget compliance by product group)
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# PRESUPUESTO (C), VENTA (D), POR CUMPLIR (E), CUMPLIMIENTO (F)
# Row 2 - 240X120 PORCELANATO
$wsCumpl.Range("C2").Value = 2826.66
$wsCumpl.Range("D2").Value = 0
$wsCumpl.Range("E2").Value = 2826.66
$wsCumpl.Range("F2").Value = 0

# Row 3 - 240X80 PORCELANATO
$wsCumpl.Range("C3").Value = 6623.26
$wsCumpl.Range("D3").Value = 1866.24
$wsCumpl.Range("E3").Value = 4757.02
$wsCumpl.Range("F3").Value = 0.2817706084314975

# Row 4 - FREGADEROS DE COCINA
$wsCumpl.Range("C4").Value = 844.7
$wsCumpl.Range("D4").Value = 0
$wsCumpl.Range("E4").Value = 844.7
$wsCumpl.Range("F4").Value = 0

# Row 5 - GRIFERIAS
$wsCumpl.Range("C5").Value = 86.41
$wsCumpl.Range("D5").Value = 0
$wsCumpl.Range("E5").Value = 86.41
$wsCumpl.Range("F5").Value = 0

# Row 6 - INODOROS
$wsCumpl.Range("C6").Value = 2907.58368146026
$wsCumpl.Range("D6").Value = 0
$wsCumpl.Range("E6").Value = 2907.58368146026
$wsCumpl.Range("F6").Value = 0

# Row 7 - LAVABOS
$wsCumpl.Range("C7").Value = 1320
$wsCumpl.Range("D7").Value = 23.4
$wsCumpl.Range("E7").Value = 1296.6
$wsCumpl.Range("F7").Value = 0.01772727272727273

# Row 8 - NO RESURTIBLES
$wsCumpl.Range("C8").Value = 415
$wsCumpl.Range("D8").Value = 0
$wsCumpl.Range("E8").Value = 415
$wsCumpl.Range("F8").Value = 0

# Row 9 - OTROS
$wsCumpl.Range("C9").Value = 0
$wsCumpl.Range("D9").Value = 0
$wsCumpl.Range("E9").Value = 0
$wsCumpl.Range("F9").Value = 0

# Row 10 - PANELES DECORATIVOS
$wsCumpl.Range("C10").Value = 4312
$wsCumpl.Range("D10").Value = 405.57
$wsCumpl.Range("E10").Value = 3906.43
$wsCumpl.Range("F10").Value = 0.0940561224489796

# Row 11 - PIEDRA SINTERIZADA
$wsCumpl.Range("C11").Value = 14235.99
$wsCumpl.Range("D11").Value = 4962.1
$wsCumpl.Range("E11").Value = 9273.889999999999
$wsCumpl.Range("F11").Value = 0.3485602336051093

# Row 12 - PORCELANATO
$wsCumpl.Range("C12").Value = 64944
$wsCumpl.Range("D12").Value = 2869.79
$wsCumpl.Range("E12").Value = 62074.21
$wsCumpl.Range("F12").Value = 0.04418868563685636

# Row 13 - PUERTAS DE SEGURIDAD (unchanged)

# Row 14 - TOTAL
$wsCumpl.Range("C14").Value = 98956.25685923838
$wsCumpl.Range("D14").Value = 10127.1
$wsCumpl.Range("E14").Value = 88829.15685923837
$wsCumpl.Range("F14").Value = 0.1023391579413258

# Column width tweaks on this sheet (D/E/F), matching Excel's re-autofit
# after the number lengths changed. ColumnWidth setter applies a fixed
# +5/6 character offset versus the raw OOXML <col width> value, so we
# compensate to land exactly on the target widths (13 / 22 / 25).
$wsCumpl.Columns.Item(4).ColumnWidth = 13 - (5/6)
$wsCumpl.Columns.Item(5).ColumnWidth = 22 - (5/6)
$wsCumpl.Columns.Item(6).ColumnWidth = 25 - (5/6)
